$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.990.17'
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").Value = '2.539.07'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '''317.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").Value = '''96.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.52%  '

$ws.Range("D7").Value = '''0.576'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Value = '''0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("D10").Value = '''35.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.59%  '

$ws.Range("D11").Value = '''0.0819'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("D12").Value = '''7.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("E13").Value = '  -4.33%  '

$ws.Range("D14").Value = '2.930.40'
$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").Value = '2.562.62'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").Value = '''15.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.46%  '

$ws.Range("D17").Value = '''0.850'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.56%  '

$ws.Range("D18").Value = '43.070.54'
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").Value = '''6.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.37%  '

$ws.Range("D20").Value = '''12.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.74%  '

$ws.Range("D21").Value = '0.0₃0967'
$ws.Range("E21").Value = '  -0.21%  '

$ws.Range("D22").Value = '''69.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.70%  '

$ws.Range("D23").Value = '''253.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").Value = '''2.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.79%  '

$ws.Range("E25").Value = '  +2.21%  '

$ws.Range("D26").Value = '''26.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.26%  '

$ws.Range("E27").Value = '  +0.13%  '

$ws.Range("E28").Value = '  +2.04%  '

$ws.Range("D29").Value = '''40.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.97%  '

$ws.Range("D30").Value = '''10.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.19%  '

$ws.Range("D31").Value = '''5.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.92%  '

$ws.Range("D32").Value = '''156.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.06%  '

$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").Value = '''19.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.44%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '''3.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.45%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.25%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''2.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.47%  '

$ws.Range("D37").Value = '''0.0800'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.18%  '

$ws.Range("E38").Value = '  +2.31%  '

$ws.Range("E39").Value = '  +1.36%  '

$ws.Range("E40").Value = '  -0.45%  '

$ws.Range("D41").Value = '''21.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.52%  '

$ws.Range("D42").Value = '''3.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.52%  '

$ws.Range("D43").Value = '''0.0304'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = '''3.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.33%  '

$ws.Range("D46").Value = '2.001.83'
$ws.Range("E46").Value = '  -1.32%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''9.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.90%  '

$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").Value = '''84.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.19%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.785.89'
$ws.Range("E49").Value = '  +0.00%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '''74.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.38%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''104.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.84%  '
